# Update automatico via Actualizar 03-10-2021 14-54-32
#
# The "Actualizar" refresh shifts the availability timestamps down one
# block (each 14-row block takes on the timestamp of the block above it)
# and stamps the newest block (rows 2:15) with the current refresh time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Oldest block (rows 30:43) takes what used to be in rows 16:29.
$ws.Range("D30:D43").Value = 44264.73777855324

# Middle block (rows 16:29) takes what used to be in rows 2:15.
$ws.Range("D16:D29").Value = 44264.75935748842

# Newest block (rows 2:15) is stamped with the current refresh timestamp.
$ws.Range("D2:D15").Value = 44265.61992098981
